$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("countries")

# Rename the "currency" header to "currency_icon" and add a new "currency_text" header
$ws.Range("C1").Value = "currency_icon"
$ws.Range("D1").Value = "currency_text"
$ws.Range("D1").Font.Bold = $true

# China keeps its icon-style value under currency_icon, clear anything in D
$ws.Range("C2").Value = "yen"
$ws.Range("D2").Value = ""

# Vietnam's currency code moves from C (currency) to D (currency_text)
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "VND"

# Cambodia's currency code moves from C (currency) to D (currency_text), and changes KH -> KHR
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "KHR"

# Match the column widths / best-fit seen in the target workbook
$ws.Columns.Item(3).ColumnWidth = 13.0
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666

# Update the selection to the new D1 cell (single-cell selection), matching target sheetView
$ws.Range("D1").Select()
